# The hyperlink wrapping "Le Monde mathematical puzzle" is removed while
# keeping the run's text and formatting (blue color, underline, etc.) intact.
# Word's Hyperlink.Delete removes only the hyperlink field/wrapper, leaving
# the display text run in place.
$d = $word.ActiveDocument

$h = $d.Hyperlinks.Item(1)
$h.Delete()
